# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted as row 14 (pushing the former
# rows 14-18 down to 15-19); every other row's data is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 14 - shifts rows 14:18 down to 15:19.
$ws.Rows("14").Insert()

# Populate the newly inserted row 14 with the new weekly record.
$ws.Cells.Item(14, 1).Value = 1
$ws.Cells.Item(14, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(14, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(14, 4).Value = 44874
$ws.Cells.Item(14, 5).Value = 15
$ws.Cells.Item(14, 6).Value = "Fruta"
$ws.Cells.Item(14, 7).Value = 100107
$ws.Cells.Item(14, 8).Value = "Otros"
$ws.Cells.Item(14, 9).Value = 100107002
$ws.Cells.Item(14, 10).Value = "Chirimoya"
$ws.Cells.Item(14, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(14, 12).Value = "Segunda"
$ws.Cells.Item(14, 13).Value = 250
$ws.Cells.Item(14, 14).Value = 22000
$ws.Cells.Item(14, 15).Value = 23000
$ws.Cells.Item(14, 16).Value = 22500
$ws.Cells.Item(14, 17).Value = "$/caja 12 kilos"
$ws.Cells.Item(14, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(14, 19).Value = 1875
$ws.Cells.Item(14, 20).Value = 12

# Keep the date column's custom date/time number format consistent with
# the rest of the "Fecha" column (style index 2 in styles.xml).
$ws.Cells.Item(14, 4).NumberFormat = $ws.Cells.Item(13, 4).NumberFormat
